# Update sample name labels in Sheet1 (column A, rows 9-14) to use the
# "TMP_" prefix, reflecting the renamed sample sources for the June event.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Range("A9").Value  = "TMP_FW_SOURCE_HR4"
$ws1.Range("A10").Value = "TMP_ESTUARY_BARGE_HR8"
$ws1.Range("A11").Value = "TMP_FW_SOURCE_HR5"
$ws1.Range("A12").Value = "TMP_FW_SOURCE_HR7"
$ws1.Range("A13").Value = "TMP_SW_SOURCE_HR7"
$ws1.Range("A14").Value = "TMP_FW_SOURCE_HR0"

# Make Sheet1 the active sheet/tab and set its selection, mirroring the
# workbook being left with Sheet1 active (instead of "Dilution sheet").
$ws1.Activate()
$ws1.Range("A18").Select()
